# Insert a new data row at row 147 (pushing existing rows 147-211 down to
# 148-212) and populate it with the new observation, matching the target
# diff: dimension grows from A1:R211 to A1:R212.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 147..211 down by one, creating a blank row 147.
$ws.Range("A147").EntireRow.Insert()

# Populate the newly inserted row 147 with the new record.
$ws.Range("A147").Value = 7
$ws.Range("B147").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C147").Value = "Ñuble"
$ws.Range("D147").Value = 44572
$ws.Range("E147").Value = 16
$ws.Range("F147").Value = 100112002
$ws.Range("G147").Value = "Pimiento"
$ws.Range("H147").Value = "Cuatro cascos verde"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 120
$ws.Range("K147").Value = 9500
$ws.Range("L147").Value = 10000
$ws.Range("M147").Value = 9750
$ws.Range("N147").Value = "$/caja 15 kilos"
$ws.Range("O147").Value = "Región del Maule"
$ws.Range("P147").Value = 650
$ws.Range("Q147").Value = 15
$ws.Range("R147").Value = "Hortaliza"
